# Update the "Datos actualizados..." timestamp string in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 20:22"

# Update numeric stats for countries whose row position does not change
# (rows above the Emiratos Arabes Unidos re-insertion point, plus the ones
# below it updated here using their pre-shift row numbers -- the later
# row-insert/delete for Emiratos Arabes Unidos will carry these new values
# down with their rows automatically, just like a real Excel edit would).

# Estados Unidos (row 4)
$ws.Range("B4").Value = 489268
$ws.Range("C4").Value = 20702
$ws.Range("D4").Value = 26187
$ws.Range("E4").Value = 445066
$ws.Range("F4").Value = 10896
$ws.Range("G4").Value = 1324
$ws.Range("H4").Value = 18015

# Canada (row 16)
$ws.Range("B16").Value = 22046
$ws.Range("C16").Value = 1281
$ws.Range("D16").Value = 5834
$ws.Range("E16").Value = 15656
$ws.Range("F16").Value = 518
$ws.Range("G16").Value = 47
$ws.Range("H16").Value = 556

# Austria (row 19)
$ws.Range("B19").Value = 13551
$ws.Range("C19").Value = 307
$ws.Range("D19").Value = 6064
$ws.Range("E19").Value = 7168
$ws.Range("F19").Value = 261
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = 319

# India (row 24)
$ws.Range("B24").Value = 7598
$ws.Range("C24").Value = 873
$ws.Range("D24").Value = 774
$ws.Range("E24").Value = 6578
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 19
$ws.Range("H24").Value = 246

# Irlanda (row 26)
$ws.Range("B26").Value = 7054
$ws.Range("C26").Value = 480
$ws.Range("D26").Value = 25
$ws.Range("E26").Value = 6742
$ws.Range("F26").Value = 194
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = 287

# Chile (row 27)
$ws.Range("B27").Value = 6501
$ws.Range("C27").Value = 529
$ws.Range("D27").Value = 1571
$ws.Range("E27").Value = 4865
$ws.Range("F27").Value = 70
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 65

# Noruega (row 28)
$ws.Range("B28").Value = 6298
$ws.Range("C28").Value = 79
$ws.Range("D28").Value = 32
$ws.Range("E28").Value = 6153
$ws.Range("F28").Value = 70
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 113

# Tunez (row 80, pre-shift)
$ws.Range("B80").Value = 671
$ws.Range("C80").Value = 28
$ws.Range("D80").Value = 25
$ws.Range("E80").Value = 621
$ws.Range("F80").Value = 78
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 25

# Liechtenstein (row 137, pre-shift)
$ws.Range("B137").Value = 79
$ws.Range("C137").Value = 1
$ws.Range("D137").Value = 55
$ws.Range("E137").Value = 23
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 1

# Belice (row 190, pre-shift)
$ws.Range("B190").Value = 10
$ws.Range("C190").Value = 1
$ws.Range("D190").Value = 0
$ws.Range("E190").Value = 8
$ws.Range("F190").Value = 1
$ws.Range("G190").Value = 1
$ws.Range("H190").Value = 2

# Re-position "Emiratos Arabes Unidos": it moves up from just below Serbia
# to just below Mexico/above Luxemburgo, with refreshed totals. Remove its
# old row (old row 44, stale data) and insert a fresh row just after Mexico
# (row 41) / before Luxemburgo (row 42), so Luxemburgo and Serbia shift back
# down into the rows Emiratos Arabes Unidos used to occupy -- their own
# figures are unchanged between the 19:52 and 20:22 snapshots.
$ws.Rows(44).Delete()
$ws.Rows(42).Insert()
$ws.Range("A42").Value = "Emiratos Arabes Unidos"
$ws.Range("B42").Value = 3360
$ws.Range("C42").Value = 370
$ws.Range("D42").Value = 418
$ws.Range("E42").Value = 2926
$ws.Range("F42").Value = 1
$ws.Range("G42").Value = 2
$ws.Range("H42").Value = 16
